# "Endpoint com exportacao de dados"
#
# Content fix on the CORPO sheet, row 14 (the "Média" / average-per-day field):
#   - Tipo (Type) column corrected from "A" (alfanumerico) to "R" (valor/real),
#     matching the other numeric/currency fields (Meta, Arrecadado, Valor doação).
#   - Descrição column typo fixed: "Media de valores poder dia"
#                               -> "Media de valores por dia"
#
# Also replays the view/selection state left behind by the edit: CORPO becomes
# the active/selected tab, and each sheet keeps the selection the author left
# on it (HEADER -> F12, CORPO -> A15:D15, TRAILER -> whole A:F columns).

$wb      = $excel.ActiveWorkbook
$header  = $wb.Worksheets.Item("HEADER")
$corpo   = $wb.Worksheets.Item("CORPO")
$trailer = $wb.Worksheets.Item("TRAILER")

# --- Data/content corrections -------------------------------------------------
$corpo.Range("C14").Value = "R"
$corpo.Range("F14").Value = "Media de valores por dia"

# --- View state (selection per sheet + active tab) ----------------------------
# Order matters: activate HEADER/TRAILER first so their selections are recorded,
# then activate CORPO last so it ends up as the workbook's active/selected sheet.
$header.Activate()
$header.Range("F12").Select()

$trailer.Activate()
$trailer.Range("A1:F1048576").Select()

$corpo.Activate()
$corpo.Range("A15:D15").Select()
